$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cells = @(
    "AA2",
    "AB4",
    "Y6",
    "AA6",
    "AB6",
    "AA7",
    "AD7",
    "Y8",
    "AA9",
    "AB9",
    "Y13",
    "AA15",
    "AB15",
    "Y17",
    "AB18",
    "Y19",
    "Y20",
    "Y23",
    "Y24",
    "Y25",
    "Y26",
    "Y30",
    "Y32",
    "AD34",
    "Y41",
    "Y43",
    "Y45",
    "Y46",
    "AD47",
    "Y48",
    "AB48",
    "Y49",
    "Y50",
    "Y51",
    "Y58",
    "AB58",
    "Y60",
    "Y62",
    "Y63",
    "Y65",
    "AA67",
    "Y69",
    "Y72",
    "Y73",
    "Y74",
    "Y75",
    "Y76",
    "Y80",
    "Y82",
    "AA84",
    "Y87",
    "AB89",
    "Y91",
    "Y92",
    "Y95",
    "AB99",
    "Y100",
    "Y101",
    "AA101",
    "Y102",
    "Y104",
    "Y105",
    "Y106",
    "Y107",
    "AD111",
    "Y112",
    "AA112",
    "AB115",
    "Y116",
    "Y117",
    "AB117",
    "AA120",
    "Y121",
    "Y123",
    "Y124",
    "AB124",
    "Y127",
    "Y128",
    "AA128",
    "AB128",
    "AD128",
    "Y129",
    "Y130",
    "AA131",
    "Y134",
    "Y136",
    "AA138",
    "AB138",
    "AB139",
    "Y145",
    "AB145",
    "AA147",
    "AB147",
    "Y150",
    "AD150",
    "Y151",
    "Y154",
    "AB157",
    "AB160",
    "Y161",
    "Y163",
    "Y165",
    "AB166",
    "AA167",
    "Y169",
    "Y171",
    "AA173",
    "Y174",
    "Y175",
    "Y180",
    "Y184",
    "Y185",
    "AA186",
    "Y188",
    "Y191",
    "Y192",
    "AA192",
    "AA193",
    "Y194",
    "AB196"
)

foreach ($cellRef in $cells) {
    $ws.Range($cellRef).Value = 0
}
